# PenjelasanDB.xlsx update:
#   - "Gedung" (building) row: INFO_5 placeholder changed from the
#     "-----" filler to a new "<<No IMB>>" (building permit number) token.
#   - Selection cursor moved from G9 to G10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 is the "Gedung" row; column G (INFO_5) currently holds the
# "-----" placeholder and needs to become "<<No IMB>>".
$target = $ws.Range("G3")
$target.Value = "<<No IMB>>"

# Writing a new value resets the cell to the plain "text" style (s=5).
# Re-apply the formatting (which also carries the quote-prefix flag that
# keeps the leading "<<" from being misread) from its still-untouched
# neighbour H3, which shares the same original look, without touching
# the text we just wrote.
$ws.Range("H3").Copy() | Out-Null
$target.PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Move the active selection to G10, matching the saved view state.
$ws.Range("G10").Select() | Out-Null
